# Apply the "remove RG to user and models certificate" edit:
#  - Shrink the body textbox (Rectangle 5) now that a line of text is gone.
#  - Rework the "portador do RG nº {{RG}} e CPF nº " sentence down to
#    "portador do CPF nº " (drop the RG reference).
#  - Merge "São Carlos" + ", " into a single run "São Carlos, ".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)   # "Rectangle 5" - the certificate body paragraph
$tr = $shape.TextFrame.TextRange

# --- 2. Replace "portador do RG nº {{RG}} e CPF nº " with "portador do CPF nº " ---
$span = $tr.Find("portador do RG n") 
$spanFull = $tr.Characters($span.Start, 34)
$spanFull.Text = "portador do CPF nº "

# Split the replacement into three runs matching the target wording breaks.
$run1 = $tr.Find("portador ")
$run1.Text = "portador "

$run2 = $tr.Find("do CPF ")
$run2.Text = "do CPF "

$run3 = $tr.Find("nº ")
$run3.Text = "nº "

# --- 3. Merge "São Carlos" and ", " into a single run "São Carlos, " ---
$sc = $tr.Find("São Carlos")
$merged = $tr.Characters($sc.Start, $sc.Length + 2)
$merged.Text = "São Carlos, "

# --- 1. Resize the shape (height shrinks because the paragraph got shorter) ---
$shape.Height = 155.1
